$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44691
$ws.Range("J2").Value = 100

# Row 3
$ws.Range("D3").Value = 44692
$ws.Range("J3").Value = 120

# Row 5
$ws.Range("D5").Value = 44687
$ws.Range("J5").Value = 160
